$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price (D) / 1h volume-change (E) figures refreshed by the
# scheduled "Updated symbol list" GitHub Actions job. Values are written
# with a leading apostrophe so Excel stores them as literal text (matching
# the existing text-formatted Price/Volume(1h) columns) rather than coercing
# them into numbers/percentages.

$ws.Range("D2").Value = "'292.40"
$ws.Range("E2").Value = "'2.03%"

$ws.Range("D3").Value = "'29.57"
$ws.Range("E3").Value = "'3.31%"

$ws.Range("D4").Value = "'5.213"
$ws.Range("E4").Value = "'2.81%"

$ws.Range("D5").Value = "'0.07135"
$ws.Range("E5").Value = "'7.34%"

$ws.Range("D6").Value = "'7.520"
$ws.Range("E6").Value = "'1.91%"

$ws.Range("D7").Value = "'3.615"
$ws.Range("E7").Value = "'6.06%"

$ws.Range("D8").Value = "'1.408"

$ws.Range("D9").Value = "'0.9114"
$ws.Range("E9").Value = "'-2.92%"

$ws.Range("D10").Value = "'0.1630"
$ws.Range("E10").Value = "'3.54%"

$ws.Range("D11").Value = "'0.07655"
$ws.Range("E11").Value = "'15.73%"

$ws.Range("D12").Value = "'0.07723"
$ws.Range("E12").Value = "'1.91%"

$ws.Range("D13").Value = "'0.02947"
$ws.Range("E13").Value = "'0.19%"

$ws.Range("D14").Value = "'0.09011"
$ws.Range("E14").Value = "'0.27%"

$ws.Range("D15").Value = "'0.001599"
$ws.Range("E15").Value = "'0.12%"

$ws.Range("D16").Value = "'0.0006564"
$ws.Range("E16").Value = "'1.31%"

$ws.Range("D17").Value = "'0.006156"
$ws.Range("E17").Value = "'-1.96%"

$ws.Range("D18").Value = "'3.489"
$ws.Range("E18").Value = "'1.32%"

$ws.Range("E19").Value = "'-0.84%"

$ws.Range("D20").Value = "'0.3270"
$ws.Range("E20").Value = "'1.69%"

$ws.Range("D21").Value = "'0.1368"
$ws.Range("E21").Value = "'5.49%"

$ws.Range("D22").Value = "'4.068"
$ws.Range("E22").Value = "'0.07%"

$ws.Range("D23").Value = "'0.1599"
$ws.Range("E23").Value = "'3.02%"

$ws.Range("D24").Value = "'0.04527"
$ws.Range("E24").Value = "'0.59%"

$ws.Range("D25").Value = "'0.001211"
$ws.Range("E25").Value = "'2.28%"

$ws.Range("D26").Value = "'0.004254"
$ws.Range("E26").Value = "'2.72%"

$ws.Range("D27").Value = "'0.0001169"
$ws.Range("E27").Value = "'-6.45%"

$ws.Range("D28").Value = "'0.0001689"
$ws.Range("E28").Value = "'4.38%"

$ws.Range("D40").Value = "'0.04385"
$ws.Range("E40").Value = "'4.36%"

$ws.Range("D41").Value = "'0.007017"
$ws.Range("E41").Value = "'3.92%"

$ws.Range("D42").Value = "'0.1279"
$ws.Range("E42").Value = "'2.38%"

$ws.Range("D43").Value = "'0.002209"
$ws.Range("E43").Value = "'9.34%"

$ws.Range("D44").Value = "'0.01352"
$ws.Range("E44").Value = "'9.71%"

$ws.Range("D45").Value = "'0.00005844"
$ws.Range("E45").Value = "'4.37%"

$ws.Range("D47").Value = "'0.01299"
$ws.Range("E47").Value = "'-0.60%"
